# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the c37a63a9-63d2-40bb-a921-0881e441ccc2 row (row 7) across the
# Overview, zh-cn and de-de sheets after a handoff xliff generation.

$wb = $excel.ActiveWorkbook
$dateFormat = "yyyy-mm-dd HH:mm:ss"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-18 10:45:51"
$overview.Range("G7").NumberFormat = $dateFormat

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-18 10:45:47"
$zhcn.Range("H7").NumberFormat = $dateFormat

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-18 10:45:51"
$dede.Range("H7").NumberFormat = $dateFormat
